# Update column G ("K") values for rows 2-23 as per the regenerated save_data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 7
    4  = 5
    5  = 3
    6  = 5
    7  = 6
    8  = 8
    9  = 5
    10 = 9
    11 = 3
    12 = 5
    13 = 9
    14 = 6
    15 = 3
    16 = 6
    17 = 12
    18 = 7
    19 = 10
    20 = 2
    21 = 4
    22 = 6
    23 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
